$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 6
$ws1.Range("F3").Value = 523
$ws1.Range("F7").Value = 183
$ws1.Range("F9").Value = 997
$ws1.Range("F10").Value = 789
$ws1.Range("F11").Value = 225
$ws1.Range("F14").Value = 803
$ws1.Range("F16").Value = 572
$ws1.Range("F17").Value = 498
$ws1.Range("F18").Value = 1317
$ws1.Range("F20").Value = 498
$ws1.Range("F21").Value = 1145
$ws1.Range("F22").Value = 2839
$ws1.Range("F23").Value = 1368
$ws1.Range("F24").Value = 679
$ws1.Range("F26").Value = 1261
$ws1.Range("F28").Value = 993
$ws1.Range("F29").Value = 343
$ws1.Range("F30").Value = 2447
$ws1.Range("C31").Value = "广州·运动番6.0-排球少年之宿命召集"
$ws1.Range("F31").Value = 458
$ws1.Range("F32").Value = 455
$ws1.Range("F33").Value = 1372

# Sheet 3: 本地生活 (Local life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 726

# Sheet 4: 全部类型 (All types - combined/mirrored view)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6
$ws4.Range("F3").Value = 726
$ws4.Range("F4").Value = 523
$ws4.Range("F13").Value = 183
$ws4.Range("F16").Value = 997
$ws4.Range("F17").Value = 789
$ws4.Range("F18").Value = 225
$ws4.Range("F26").Value = 803
$ws4.Range("F28").Value = 572
$ws4.Range("F29").Value = 498
$ws4.Range("F30").Value = 1317
$ws4.Range("F32").Value = 498
$ws4.Range("F33").Value = 1145
$ws4.Range("F34").Value = 2839
$ws4.Range("F35").Value = 1368
$ws4.Range("F36").Value = 679
$ws4.Range("F38").Value = 1261
$ws4.Range("F42").Value = 993
$ws4.Range("F43").Value = 343
$ws4.Range("F44").Value = 2447
$ws4.Range("C45").Value = "广州·运动番6.0-排球少年之宿命召集"
$ws4.Range("F45").Value = 458
$ws4.Range("F46").Value = 455
$ws4.Range("F47").Value = 1372
